$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataTypes")
$ws.Activate()

# Row 4 corresponds to the "char" data type. Update the MySQL (C4) and
# SQLite (D4) template cells from "CHAR" to "CHAR(0)" (sized template,
# matching the new "mysql"/"sqlite t4" templates).
$ws.Range("C4").Value = "CHAR(0)"
$ws.Range("D4").Value = "CHAR(0)"

# Update the remembered selection on the DataTypes sheet.
$ws.Range("D5").Select()
